$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the customer transaction rows (2 through 4), including the
# "Order #" column, leaving the cells empty instead of reordering/skipping
# rows when an exception (missing/blank order) is hit.
$ws.Range("B2:H4").ClearContents()
$ws.Range("A3:A4").ClearContents()
